# Applies the "Corrected excel sheets for application fix issues" edit.
#
# Summary of changes:
#  1. Summary sheet       - shrink the stale selection from A7:XFD15 to A7:XFD14.
#  2. Repayment schedule  - add a new "O" column (rows 2-15) mirroring the
#                            existing "N" column's values/formatting.
#  3. Transactions        - bump the running transaction IDs in A2:A4 by +28
#                            (87->115, 86->114, 85->113) and move the active
#                            selection to the single cell D3.

$wb = $excel.ActiveWorkbook

# --- 1. Summary sheet: fix selection -----------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate()
$wsSummary.Range("A7:XFD14").Select()

# --- 2. Repayment schedule: insert column O, rows 2-15 ------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Range("N2:N15").Copy($wsRepay.Range("O2:O15"))

# --- 3. Transactions: update IDs + selection -----------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("A2").Value = 115
$wsTrans.Range("A3").Value = 114
$wsTrans.Range("A4").Value = 113
$wsTrans.Activate()
$wsTrans.Range("D3").Select()
